$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = -1
    3  = 6
    4  = 1
    5  = 1
    6  = -1
    7  = 3
    8  = -1
    9  = -3
    12 = 1
    13 = -4
    14 = 0
    15 = 6
    16 = 2
    17 = 0
    18 = -5
    19 = 4
    21 = 1
    22 = 1
    23 = -3
    24 = 2
    26 = -1
    27 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
